# Report Generator Via Outlook
# Rewrites the Test_Execution_Report data grid (Sheet1!A2:F20) with the
# latest test-run results: re-numbered/renumbered test cases, refreshed
# Source/Destination counts, updated Validation Status, and four new rows
# (17-20) for the newly added SMA2 / NPA portfolio-health checks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2,  "Test Case 2",  251496,             251499,     "FAIL", "Bank_Overview_Widget",               "Customers_count"),
    @(3,  "Test Case 3",  251496,             251499,     "FAIL", "Bank_Overview_Widget",               "Active_Loans"),
    @(4,  "Test Case 4",  4539016406.480295,  4539092685, "PASS", "Product_Overview",                    "Total Outstanding"),
    @(5,  "Test Case 8",  1051,               1046,       "FAIL", "Portfolio_Loan_Disbursed",            "Customers"),
    @(6,  "Test Case 9",  1051,               1046,       "FAIL", "Portfolio_Loan_Disbursed",            "Loans"),
    @(7,  "Test Case 10", 54583000,           50373000,   "FAIL", "Portfolio_Loan_Disbursed_Today",      "Loan disbursed"),
    @(8,  "Test Case 12", 1051,               1046,       "FAIL", "Portfolio_Loan_Disbursed_Today",      "Loans"),
    @(9,  "Test Case 13", 54583000,           50373000,   "FAIL", "Portfolio_Loan_Disbursed_Yesterday",  "Loan disbursed"),
    @(10, "Test Case 15", 1051,               1046,       "FAIL", "Portfolio_Loan_Disbursed_Yesterday",  "Loans"),
    @(11, "Test Case 17", 4539016406.480296,  4539092685, "PASS", "Portfolio_Items_Widget",              "OutStanding"),
    @(12, "Test Case 19", 1837928000,         1802201000, "FAIL", "Portfolio_Items_Widget",              "Loans Disbursed"),
    @(13, "Test Case 20", 3317057243.480174,  3402313794, "PASS", "Portfolio_Health_Category_OTR",       "Total Outstanding"),
    @(14, "Test Case 21", 188014,             189809,     "FAIL", "Portfolio_Health_Category_OTR",       "Total Customers"),
    @(15, "Test Case 23", 14716,              14717,      "FAIL", "Portfolio_Health_Category_SMA0",      "Total Customers"),
    @(16, "Test Case 25", 11246,              11246,      "PASS", "Portfolio_Health_Category_SMA1",      "Total Customers"),
    @(17, "Test Case 26", 106903206.86,       106903299,  "FAIL", "Portfolio_Health_Category_SMA2",      "Total Outstanding"),
    @(18, "Test Case 27", 6637,               6637,       "PASS", "Portfolio_Health_Category_SMA2",      "Total Customers"),
    @(19, "Test Case 28", 720633169.6399996,  720633275,  "FAIL", "Portfolio_Health_Category_NPA",       "Total Outstanding"),
    @(20, "Test Case 29", 30883,              30883,      "PASS", "Portfolio_Health_Category_NPA",       "Total Customers")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
    $ws.Cells.Item($rowNum, 6).Value = $r[6]
}
